$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Simple field updates on the "Metadata" sheet ---
$ws1.Range("B3").Value = "0.1.7"
$ws1.Range("B6").Value = "draft"
$ws1.Range("B8").Value = "2024-08-23T10:17:11-05:00"
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Insert a new "Jurisdiction" row after row 11 (Contact), pushing
#     Description/Purpose/Copyright/Immutable down by one row ---
# Shift rows 12:15 down to 13:16 (copy including formatting) first so the
# original formatting (borders/fill/alignment) of every row is preserved
# without creating brand-new, unused cell styles.
$ws1.Rows("12:15").Copy()
$ws1.Rows("13:16").PasteSpecial(-4104)

# Now set row 12 to the new Jurisdiction field (keeps the style/format that
# was already on that row - i.e. the same as the other data rows).
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""
